$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2166.6667
$ws.Range("I62").Value = 2498.75
$ws.Range("J62").Value = 1502.5
$ws.Range("K62").Value = 2498.75
$ws.Range("L62").Value = 1502.5
$ws.Range("M62").Value = -1874.75
$ws.Range("N62").Value = -2750.5

$ws.Range("H64").Value = 2864.7273
$ws.Range("J64").Value = 2901
$ws.Range("L64").Value = 2901
$ws.Range("N64").Value = -3397

$ws.Range("H65").Value = 2166.6667
$ws.Range("I65").Value = 2498.75
$ws.Range("J65").Value = 1502.5
$ws.Range("K65").Value = 12493.75
$ws.Range("L65").Value = 7512.5
$ws.Range("M65").Value = -9373.75
$ws.Range("N65").Value = -13752.5

$ws.Range("H67").Value = 2864.7273
$ws.Range("J67").Value = 2901
$ws.Range("L67").Value = 2901
$ws.Range("N67").Value = -4617

$ws.Range("H70").Value = 2500749.8
$ws.Range("I70").Value = 999.5
$ws.Range("J70").Value = 5000500
$ws.Range("K70").Value = 2998.5
$ws.Range("L70").Value = 15001500
$ws.Range("M70").Value = -2728.5
$ws.Range("N70").Value = -15002040

$ws.Range("H73").Value = 2500749.8
$ws.Range("I73").Value = 999.5
$ws.Range("J73").Value = 5000500
$ws.Range("K73").Value = 2998.5
$ws.Range("L73").Value = 15001500
$ws.Range("M73").Value = -2062.5
$ws.Range("N73").Value = -15003372

$ws.Range("H137").Value = 1842.4546
$ws.Range("I137").Value = 1424.15
$ws.Range("J137").Value = 2486
$ws.Range("K137").Value = 4272.450000000001
$ws.Range("L137").Value = 7458
$ws.Range("M137").Value = -1722.450000000001
$ws.Range("N137").Value = -12558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws.Range("H32").Value = 5386.172
$ws.Range("I32").Value = 4147.712
$ws.Range("K32").Value = 4147.712
$ws.Range("M32").Value = -3860.712

$ws.Range("H74").Value = 1327.3334
$ws.Range("I74").Value = 1037.5454
$ws.Range("J74").Value = 2124.25
$ws.Range("K74").Value = 1037.5454
$ws.Range("L74").Value = 2124.25
$ws.Range("M74").Value = -163.5454
$ws.Range("N74").Value = -3872.25

$ws.Range("H77").Value = 1327.3334
$ws.Range("I77").Value = 1037.5454
$ws.Range("J77").Value = 2124.25
$ws.Range("K77").Value = 5187.727
$ws.Range("L77").Value = 10621.25
$ws.Range("M77").Value = -819.7269999999999
$ws.Range("N77").Value = -19357.25

$ws.Range("H97").Value = 292.91666
$ws.Range("J97").Value = 414.33334
$ws.Range("L97").Value = 414.33334
$ws.Range("N97").Value = -1406.33334

$ws.Range("I110").Value = 1200
$ws.Range("K110").Value = 1200
$ws.Range("M110").Value = 845

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""

$ws.Range("H122").Value = 5784.3335
$ws.Range("I122").Value = 7448.75
$ws.Range("K122").Value = 22346.25
$ws.Range("M122").Value = -19896.25

$ws.Range("H132").Value = 3068.5908
$ws.Range("I132").Value = 2367.5334
$ws.Range("K132").Value = 7102.600199999999
$ws.Range("M132").Value = -4572.600199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""

$ws.Range("H22").Value = 749.4
$ws.Range("I22").Value = 762
$ws.Range("K22").Value = 762
$ws.Range("M22").Value = -589

$ws.Range("H99").Value = 789.1667
$ws.Range("I99").Value = 748.6
$ws.Range("K99").Value = 748.6
$ws.Range("M99").Value = 749.4

$ws.Range("H134").Value = 2874.7856
$ws.Range("I134").Value = 2807.4443
$ws.Range("K134").Value = 8422.332900000001
$ws.Range("M134").Value = -5887.332900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1325
$ws.Range("I62").Value = 1325
$ws.Range("K62").Value = 1325
$ws.Range("M62").Value = -701

$ws.Range("H65").Value = 1325
$ws.Range("I65").Value = 1325
$ws.Range("K65").Value = 6625
$ws.Range("M65").Value = -3505

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""

$ws.Range("H99").Value = 8087.1665
$ws.Range("I99").Value = 8503
$ws.Range("J99").Value = 7671.3335
$ws.Range("K99").Value = 8503
$ws.Range("L99").Value = 7671.3335
$ws.Range("M99").Value = -7005
$ws.Range("N99").Value = -10667.3335

$ws.Range("H126").Value = 8087.1665
$ws.Range("I126").Value = 8503
$ws.Range("J126").Value = 7671.3335
$ws.Range("K126").Value = 25509
$ws.Range("L126").Value = 23014.0005
$ws.Range("M126").Value = -23039
$ws.Range("N126").Value = -27954.0005

$ws.Range("H132").Value = 4415
$ws.Range("I132").Value = 4249.1113
$ws.Range("J132").Value = 4628.2856
$ws.Range("K132").Value = 12747.3339
$ws.Range("L132").Value = 13884.8568
$ws.Range("M132").Value = -10217.3339
$ws.Range("N132").Value = -18944.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2552.6667
$ws.Range("I18").Value = 594.8
$ws.Range("K18").Value = 1784.4
$ws.Range("M18").Value = -1615.4

$ws.Range("H130").Value = 7747.5
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2118
$ws.Range("I122").Value = 989.4
$ws.Range("K122").Value = 2968.2
$ws.Range("M122").Value = -518.1999999999998

$ws.Range("H126").Value = 999.5
$ws.Range("J126").Value = 999.5
$ws.Range("L126").Value = 2998.5
$ws.Range("N126").Value = -7938.5

$ws.Range("H134").Value = 73074.664
$ws.Range("J134").Value = 73074.664
$ws.Range("L134").Value = 219223.992
$ws.Range("N134").Value = -224293.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4591.7144
$ws.Range("I132").Value = 4028.8
$ws.Range("K132").Value = 12086.4
$ws.Range("M132").Value = -9556.400000000001

$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""

Write-Output "edits applied"
